$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.177.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.903.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +2.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3767"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07250"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9003"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08465"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.893.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008626"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.207.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.067"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.137.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.423"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("E26").Value = "  +4.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.751"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.26%  "
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.819"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.895"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09262"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8069"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05067"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.238"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.436"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.947"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.624"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5734"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01990"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.012"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.640"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "116.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1512"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4867"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9998"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.612"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.42%  "